$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Major User View")
$ws.Range("D42").Value = [char]0x2713
